# adding full run for ZEV Jan R2-4 and modifying files for consistency in R2-4
#
# Rows 2-22 have their roboticRNAPrep (H), RIBOSOMAL_BAND (I) and
# SMALL_RNA_BANDS (K) columns converted from text ("no"/"Y") to native
# Excel boolean values (FALSE/TRUE). The RIBOSOMAL_BAND and
# SMALL_RNA_BANDS columns also pick up a dedicated 10pt black Arial font
# once they hold TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 22; $r++) {
    $ws.Range("H$r").Value = $false

    $ws.Range("I$r").Font.Name = "Arial"
    $ws.Range("I$r").Font.Size = 10
    $ws.Range("I$r").Font.Color = 0
    $ws.Range("I$r").Value = $true

    $ws.Range("K$r").Font.Name = "Arial"
    $ws.Range("K$r").Font.Size = 10
    $ws.Range("K$r").Font.Color = 0
    $ws.Range("K$r").Value = $true
}

# Move the active selection from G34 to G30.
$ws.Range("G30").Select() | Out-Null
